$d = $word.ActiveDocument

function Replace-ExactText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Table of timing / metric results: numeric value updates ---
Replace-ExactText "1.77s" "0.54s"
Replace-ExactText "18.85s" "12.70s"
Replace-ExactText "25.27s" "15.12s"
Replace-ExactText "43554" "39588"
Replace-ExactText "4.22s" "1.76s"
Replace-ExactText "243.88s" "174.02s"
Replace-ExactText "504.51s" "192.41s"
Replace-ExactText "5.48s" "4.66s"
Replace-ExactText "535.87s" "306.87s"
Replace-ExactText "1297.20s" "941.20s"
Replace-ExactText "1006718" "820398"
Replace-ExactText "152" "153"

# --- Commentary cell text update ---
Replace-ExactText "by the way" "a bad movie"

Write-Output "done"
